$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4696.4
$ws.Range("I69").Value = 2483
$ws.Range("J69").Value = 5645
$ws.Range("K69").Value = 7449
$ws.Range("L69").Value = 16935
$ws.Range("M69").Value = -6575
$ws.Range("N69").Value = -18683
$ws.Range("H72").Value = 4696.4
$ws.Range("I72").Value = 2483
$ws.Range("J72").Value = 5645
$ws.Range("K72").Value = 22347
$ws.Range("L72").Value = 50805
$ws.Range("M72").Value = -17979
$ws.Range("N72").Value = -59541
$ws.Range("H99").Value = 2627.4736
$ws.Range("I99").Value = 734.8570999999999
$ws.Range("K99").Value = 2204.5713
$ws.Range("M99").Value = -706.5712999999996
$ws.Range("H132").Value = 48040.906
$ws.Range("I132").Value = 27186.41
$ws.Range("K132").Value = 81559.23
$ws.Range("M132").Value = -79029.23
$ws.Range("H137").Value = 1906230.4
$ws.Range("J137").Value = 5961618
$ws.Range("L137").Value = 17884854
$ws.Range("N137").Value = -17889954
$ws.Range("H138").Value = 4446.535
$ws.Range("I138").Value = 2424.182
$ws.Range("J138").Value = 4699.3296
$ws.Range("K138").Value = 7272.545999999999
$ws.Range("L138").Value = 14097.9888
$ws.Range("M138").Value = -2132.545999999999
$ws.Range("N138").Value = -24377.9888

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1350.1904
$ws.Range("I2").Value = 1350.1904
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1350.1904
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1237.1904
$ws.Range("H32").Value = 3783.68
$ws.Range("I32").Value = 2969.6191
$ws.Range("J32").Value = 8057.5
$ws.Range("K32").Value = 2969.6191
$ws.Range("L32").Value = 8057.5
$ws.Range("M32").Value = -2682.6191
$ws.Range("N32").Value = -8631.5
$ws.Range("H61").Value = 3903.3
$ws.Range("I61").Value = 3377.375
$ws.Range("K61").Value = 3377.375
$ws.Range("M61").Value = -3165.375
$ws.Range("H74").Value = 22753858
$ws.Range("I74").Value = 30119.857
$ws.Range("K74").Value = 30119.857
$ws.Range("M74").Value = -29245.857
$ws.Range("H77").Value = 22753858
$ws.Range("I77").Value = 30119.857
$ws.Range("K77").Value = 150599.285
$ws.Range("M77").Value = -146231.285
$ws.Range("H97").Value = 3520.125
$ws.Range("I97").Value = 2594.4285
$ws.Range("K97").Value = 2594.4285
$ws.Range("M97").Value = -2098.4285
$ws.Range("H110").Value = 7504.8613
$ws.Range("I110").Value = 8280.280000000001
$ws.Range("K110").Value = 8280.280000000001
$ws.Range("M110").Value = -6235.280000000001
$ws.Range("H116").Value = 1350.1904
$ws.Range("I116").Value = 1350.1904
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1350.1904
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 943.8096
$ws.Range("H136").Value = 3903.3
$ws.Range("I136").Value = 3377.375
$ws.Range("K136").Value = 10132.125
$ws.Range("M136").Value = -7582.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1350.1904
$ws.Range("I3").Value = 1350.1904
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1350.1904
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1236.1904
$ws.Range("H94").Value = 3594.2856
$ws.Range("I94").Value = 3832
$ws.Range("K94").Value = 3832
$ws.Range("M94").Value = -3381
$ws.Range("H99").Value = 5375.615
$ws.Range("I99").Value = 988.4
$ws.Range("J99").Value = 19999.666
$ws.Range("K99").Value = 988.4
$ws.Range("L99").Value = 19999.666
$ws.Range("M99").Value = 509.6
$ws.Range("N99").Value = -22995.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3113
$ws.Range("I16").Value = 3138.182
$ws.Range("J16").Value = 2974.5
$ws.Range("K16").Value = 3138.182
$ws.Range("L16").Value = 2974.5
$ws.Range("M16").Value = -2851.182
$ws.Range("N16").Value = -3548.5
$ws.Range("H31").Value = 22537032
$ws.Range("I31").Value = 3839608.8
$ws.Range("K31").Value = 3839608.8
$ws.Range("M31").Value = -3839313.8
$ws.Range("H34").Value = 22537032
$ws.Range("I34").Value = 3839608.8
$ws.Range("K34").Value = 3839608.8
$ws.Range("M34").Value = -3839406.8
$ws.Range("H39").Value = 4124.0835
$ws.Range("I39").Value = 3248.1667
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 3248.1667
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -2857.1667
$ws.Range("N39").Value = -5782
$ws.Range("H49").Value = 4124.0835
$ws.Range("I49").Value = 3248.1667
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 3248.1667
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = -3066.1667
$ws.Range("N49").Value = -5364
$ws.Range("H99").Value = 3819.6897
$ws.Range("I99").Value = 3956.28
$ws.Range("J99").Value = 2966
$ws.Range("K99").Value = 3956.28
$ws.Range("L99").Value = 2966
$ws.Range("M99").Value = -2458.28
$ws.Range("N99").Value = -5962
$ws.Range("H113").Value = 3113
$ws.Range("I113").Value = 3138.182
$ws.Range("J113").Value = 2974.5
$ws.Range("K113").Value = 3138.182
$ws.Range("L113").Value = 2974.5
$ws.Range("M113").Value = -968.1819999999998
$ws.Range("N113").Value = -7314.5
$ws.Range("H122").Value = 1389.2858
$ws.Range("I122").Value = 1172
$ws.Range("K122").Value = 3516
$ws.Range("M122").Value = -1066
$ws.Range("H126").Value = 3819.6897
$ws.Range("I126").Value = 3956.28
$ws.Range("J126").Value = 2966
$ws.Range("K126").Value = 11868.84
$ws.Range("L126").Value = 8898
$ws.Range("M126").Value = -9398.84
$ws.Range("N126").Value = -13838
$ws.Range("H134").Value = 37351.25
$ws.Range("I134").Value = 41355.46
$ws.Range("J134").Value = 19999.666
$ws.Range("K134").Value = 124066.38
$ws.Range("L134").Value = 59998.99800000001
$ws.Range("M134").Value = -121531.38
$ws.Range("N134").Value = -65068.99800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 33561736
$ws.Range("I4").Value = 38812990
$ws.Range("K4").Value = 116438970
$ws.Range("M4").Value = -116438858
$ws.Range("H86").Value = 326.5
$ws.Range("I86").Value = 61
$ws.Range("K86").Value = 183
$ws.Range("M86").Value = 1003
$ws.Range("H89").Value = 326.5
$ws.Range("I89").Value = 61
$ws.Range("K89").Value = 549
$ws.Range("M89").Value = 5379
$ws.Range("H92").Value = 1020
$ws.Range("J92").Value = 1212.5
$ws.Range("L92").Value = 3637.5
$ws.Range("N92").Value = -6133.5
$ws.Range("H113").Value = 413.45834
$ws.Range("I113").Value = 355.88235
$ws.Range("J113").Value = 553.2857
$ws.Range("K113").Value = 1067.64705
$ws.Range("L113").Value = 1659.8571
$ws.Range("M113").Value = 1102.35295
$ws.Range("N113").Value = -5999.8571
$ws.Range("H132").Value = 1248.4572
$ws.Range("I132").Value = 1161.3478
$ws.Range("J132").Value = 1415.4166
$ws.Range("K132").Value = 10452.1302
$ws.Range("L132").Value = 12738.7494
$ws.Range("M132").Value = -7922.1302
$ws.Range("N132").Value = -17798.7494

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10132
$ws.Range("I44").Value = 11514
$ws.Range("J44").Value = 8750
$ws.Range("K44").Value = 11514
$ws.Range("L44").Value = 8750
$ws.Range("M44").Value = -10918
$ws.Range("N44").Value = -9942
$ws.Range("H113").Value = 4640.048
$ws.Range("J113").Value = 5109.778
$ws.Range("L113").Value = 5109.778
$ws.Range("N113").Value = -9449.778
$ws.Range("H114").Value = 28961
$ws.Range("J114").Value = 28961
$ws.Range("L114").Value = 28961
$ws.Range("N114").Value = -37639
$ws.Range("H132").Value = 12056.34
$ws.Range("I132").Value = 13916.952
$ws.Range("K132").Value = 41750.856
$ws.Range("M132").Value = -39220.856

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3217.9395
$ws.Range("J16").Value = 2800
$ws.Range("L16").Value = 2800
$ws.Range("N16").Value = -3140
$ws.Range("H61").Value = 5357.7144
$ws.Range("I61").Value = 4871.4546
$ws.Range("J61").Value = 7140.6665
$ws.Range("K61").Value = 4871.4546
$ws.Range("L61").Value = 7140.6665
$ws.Range("M61").Value = -4669.4546
$ws.Range("N61").Value = -7544.6665
$ws.Range("H113").Value = 5357.7144
$ws.Range("I113").Value = 4871.4546
$ws.Range("J113").Value = 7140.6665
$ws.Range("K113").Value = 4871.4546
$ws.Range("L113").Value = 7140.6665
$ws.Range("M113").Value = -2701.4546
$ws.Range("N113").Value = -11480.6665
$ws.Range("H122").Value = 406676.6
$ws.Range("I122").Value = 773815
$ws.Range("K122").Value = 2321445
$ws.Range("M122").Value = -2318995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 45000000
$ws.Range("I15").Value = 45000000
$ws.Range("K15").Value = 45000000
$ws.Range("M15").Value = -44999712
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H49").Value = 22484.5
$ws.Range("J49").Value = 22484.5
$ws.Range("L49").Value = 22484.5
$ws.Range("N49").Value = -22944.5
$ws.Range("H113").Value = 1882.7778
$ws.Range("I113").Value = 1698.5714
$ws.Range("K113").Value = 5095.7142
$ws.Range("M113").Value = -2925.7142
$ws.Range("H126").Value = 5156
$ws.Range("I126").Value = 5317.972
$ws.Range("J126").Value = 3698.25
$ws.Range("K126").Value = 15953.916
$ws.Range("L126").Value = 11094.75
$ws.Range("M126").Value = -13483.916
$ws.Range("N126").Value = -16034.75
